$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert row for "fetch" service-discovery entry (new row 5) ---
$ws.Rows("5:5").Insert()
$ws.Range("A4:D4").Copy()
$ws.Range("A5:D5").PasteSpecial(-4122)
$ws.Range("A5").Value = "fetch"
$ws.Range("B5").Value = "fetch-eureka-server"
$ws.Range("C5").Value = 8000
$ws.Range("D5").Value = ""

# --- Insert row for "reports" service-discovery entry (new row 8) ---
$ws.Rows("8:8").Insert()
$ws.Range("A7:D7").Copy()
$ws.Range("A8:D8").PasteSpecial(-4122)
$ws.Range("A8").Value = "reports"
$ws.Range("B8").Value = "reports-eureka-server"
$ws.Range("C8").Value = 8001
$ws.Range("D8").Value = ""

# --- Insert row for "persistence" service-discovery entry (new row 10) ---
$ws.Rows("10:10").Insert()
$ws.Range("A9:D9").Copy()
$ws.Range("A10:D10").PasteSpecial(-4122)
$ws.Range("A10").Value = "persistence"
$ws.Range("B10").Value = "persistence-eureka-server"
$ws.Range("C10").Value = 8002
$ws.Range("D10").Value = ""

# --- Column B needs to widen to fit the new longer service names ---
$ws.Columns("B:B").ColumnWidth = 21.8333333333333

# --- Selection / active cell matches post-edit state ---
$ws.Range("D15").Select() | Out-Null
